# Actualización automática 2025-11-07 08:30:07
#
# Records new "noviembre" sales figures for several clients of the asesor
# GUERRERO FAREZ FABIAN MAURICIO, and rolls those new figures up into the
# per-group summary sheet ("VENTAS POR GRUPO"), the "de 54" non-zero-count
# row on that same sheet, and the grand-total row on "VENTA MENSUAL".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client, per-product-group sales values.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M8").Value  = 314.8    # ARMIJOS AMBROSI LUIS KLEBER      - PORCELANATO
$wsGrupo.Range("L9").Value  = 443.44   # ASTUDILLO ESPINOZA JOSE MANUEL   - PIEDRA SINTERIZADA
$wsGrupo.Range("M12").Value = 1727.91  # BUSTAMANTE ROSERO MARCO TULIO    - PORCELANATO
$wsGrupo.Range("M36").Value = 33.7     # ORTEGA ROMAN KLEBER ERWIN        - PORCELANATO
$wsGrupo.Range("M39").Value = 121.31   # PACHECO NIVICELA SANDRA ELISABETH- PORCELANATO
$wsGrupo.Range("D48").Value = 475.2    # RUIZ TINIZARAY YOHANNA MARYURI   - 240X80 PORCELANATO
$wsGrupo.Range("I53").Value = 23.4     # WONG SANCHEZ CLAUDIA PAULINA     - LAVABOS

# Row 56 holds "<n> de 54" counts of non-zero entries per column; each
# column that received a new non-zero value above gains one count.
$wsGrupo.Range("D56").Value = "3 de 54"
$wsGrupo.Range("I56").Value = "1 de 54"
$wsGrupo.Range("L56").Value = "3 de 54"
$wsGrupo.Range("M56").Value = "5 de 54"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same clients' "noviembre" monthly totals.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F8").Value  = 314.8    # ARMIJOS AMBROSI LUIS KLEBER
$wsMensual.Range("F9").Value  = 443.44   # ASTUDILLO ESPINOZA JOSE MANUEL
$wsMensual.Range("F12").Value = 1727.91  # BUSTAMANTE ROSERO MARCO TULIO
$wsMensual.Range("F36").Value = 33.7     # ORTEGA ROMAN KLEBER ERWIN
$wsMensual.Range("F39").Value = 121.31   # PACHECO NIVICELA SANDRA ELISABETH
$wsMensual.Range("F48").Value = 475.2    # RUIZ TINIZARAY YOHANNA MARYURI
$wsMensual.Range("F55").Value = 23.4     # WONG SANCHEZ CLAUDIA PAULINA (1st line)
$wsMensual.Range("F56").Value = 23.4     # WONG SANCHEZ CLAUDIA PAULINA (2nd line)

# Grand-total row for the "noviembre" column.
$wsMensual.Range("F60").Value = 9930.85
